$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column (H) to the right of the existing "sum" column (G).
# Copy G1's formatting (bold, centered, bordered header style) onto H1 so the
# new header cell reuses the same style as the rest of the header row.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Populate the new data column with zeros (unstyled, like the other data cells).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
